$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value2 = 0.07867893939452165
$ws.Range("K2").Value2 = -0.4579186939869491
$ws.Range("G3").Value2 = 0.2834097224933642
$ws.Range("K3").Value2 = 4.954820271343294
$ws.Range("G4").Value2 = 0.6100384473891095
$ws.Range("K4").Value2 = 7.620448435112452
$ws.Range("G5").Value2 = 1.053359069074734
$ws.Range("K5").Value2 = 7.637930026244192
$ws.Range("G6").Value2 = 1.607051367755741
$ws.Range("K6").Value2 = 4.192563937751054
$ws.Range("G7").Value2 = 2.263696227475027
$ws.Range("K7").Value2 = 5.127828414156189
$ws.Range("G8").Value2 = 3.01485362797923
$ws.Range("K8").Value2 = 3.097615694988355
$ws.Range("G9").Value2 = 3.851180243942212
$ws.Range("K9").Value2 = 5.899531749309721
$ws.Range("G10").Value2 = 4.762574860949703
$ws.Range("K10").Value2 = 14.88181664455112
$ws.Range("G11").Value2 = 5.738338041100786
$ws.Range("K11").Value2 = 20.08590638391232
$ws.Range("G12").Value2 = 6.767340374005929
$ws.Range("K12").Value2 = 15.89477455061515
$ws.Range("G13").Value2 = 7.838185056218411
$ws.Range("K13").Value2 = 19.66235674573823
$ws.Range("G14").Value2 = 8.939369493888595
$ws.Range("K14").Value2 = 15.00477806561599
$ws.Range("G15").Value2 = 10.0594354568538
$ws.Range("K15").Value2 = 15.86499153445503
$ws.Range("G16").Value2 = 11.18710685144444
$ws.Range("K16").Value2 = 14.27460199283323
$ws.Range("G17").Value2 = 12.3114161914133
$ws.Range("K17").Value2 = 16.46757759209895
$ws.Range("G18").Value2 = 13.42181709598683
$ws.Range("K18").Value2 = 16.14934339717764
$ws.Range("G19").Value2 = 14.50828356327137
$ws.Range("K19").Value2 = 16.56652264531595
$ws.Range("G20").Value2 = 15.56139371975123
$ws.Range("K20").Value2 = 10.17974947234659
$ws.Range("G21").Value2 = 16.57240536301691
$ws.Range("K21").Value2 = 8.717725023027125
$ws.Range("G22").Value2 = 17.53331467743018
$ws.Range("K22").Value2 = 16.28090898770595
$ws.Range("G23").Value2 = 18.43690513843573
$ws.Range("K23").Value2 = 16.61498536166046
$ws.Range("G24").Value2 = 19.27678848726305
$ws.Range("K24").Value2 = 18.19758017687391
$ws.Range("G25").Value2 = 20.04743223365554
$ws.Range("K25").Value2 = 18.21880266112855
$ws.Range("G26").Value2 = 20.74418055680502
$ws.Range("K26").Value2 = 18.21455970610711
$ws.Range("G27").Value2 = 21.36326783575629
$ws.Range("K27").Value2 = 17.37585455779326
$ws.Range("G28").Value2 = 21.90182322228899
$ws.Range("K28").Value2 = 20.14396561101301
$ws.Range("G29").Value2 = 22.35787096932565
$ws.Range("K29").Value2 = 16.28028710310328
$ws.Range("G30").Value2 = 22.73032595258459
$ws.Range("K30").Value2 = 16.25407199910463
$ws.Range("G31").Value2 = 23.01898193551957
$ws.Range("K31").Value2 = 16.12766064845548
$ws.Range("G32").Value2 = 23.22450110330703
$ws.Range("K32").Value2 = 17.78976053630376
$ws.Range("G33").Value2 = 23.34839834846756
$ws.Range("K33").Value2 = 21.41039554918826
$ws.Range("G34").Value2 = 23.3930243426538
$ws.Range("K34").Value2 = 23.48104384516746
$ws.Range("G35").Value2 = 23.36155418262033
$ws.Range("K35").Value2 = 26.93362624950406
$ws.Range("G36").Value2 = 23.25797045704835
$ws.Range("K36").Value2 = 23.86304169882229
$ws.Range("G37").Value2 = 23.08705582812599
$ws.Range("K37").Value2 = 22.91185492742067
$ws.Range("G38").Value2 = 22.85438830034078
$ws.Range("K38").Value2 = 20.51248442757407
$ws.Range("G39").Value2 = 22.56634638413218
$ws.Range("K39").Value2 = 19.332310817687
$ws.Range("G40").Value2 = 22.23012006344638
$ws.Range("K40").Value2 = 21.30219149580196
$ws.Range("G41").Value2 = 21.85373722402712
$ws.Range("K41").Value2 = 20.66156169859631
$ws.Range("G42").Value2 = 21.44610387769768
$ws.Range("K42").Value2 = 16.76131801877803
$ws.Range("G43").Value2 = 21.01706271973296
$ws.Range("K43").Value2 = 21.98307891714498
$ws.Range("G44").Value2 = 20.57746698417667
$ws.Range("K44").Value2 = 17.99428544384139
$ws.Range("G45").Value2 = 20.13928514792284
$ws.Range("K45").Value2 = 18.53644935363875
$ws.Range("G46").Value2 = 19.71572056659708
$ws.Range("K46").Value2 = 15.50119159484117
$ws.Range("G47").Value2 = 19.32136718589418
$ws.Range("K47").Value2 = 13.25876754002618
$ws.Range("G48").Value2 = 18.9723857526461
$ws.Range("K48").Value2 = 12.76612072013405
$ws.Range("G49").Value2 = 18.6867157605099
$ws.Range("K49").Value2 = 12.1848734216308
$ws.Range("G50").Value2 = 18.48430975228594
$ws.Range("K50").Value2 = 15.98540036556252
$ws.Range("G51").Value2 = 18.38739757509961
$ws.Range("K51").Value2 = 16.59243924612798
$ws.Range("G52").Value2 = 18.42077588820602
$ws.Range("K52").Value2 = 15.47025629521547
$ws.Range("G53").Value2 = 18.61211176586196
$ws.Range("K53").Value2 = 18.56363023194041
$ws.Range("G54").Value2 = 18.99226673737393
$ws.Range("K54").Value2 = 19.47056594718946
$ws.Range("G55").Value2 = 19.59561757915424
$ws.Range("K55").Value2 = 21.23398978014418
$ws.Range("G56").Value2 = 20.46038267041456
$ws.Range("K56").Value2 = 24.26027232667296
$ws.Range("G57").Value2 = 21.62892370321029
$ws.Range("K57").Value2 = 26.90085600578941
$ws.Range("G58").Value2 = 23.14802977229899
$ws.Range("K58").Value2 = 31.11250904678676
$ws.Range("G59").Value2 = 25.06915655309442
$ws.Range("K59").Value2 = 32.6705964745689
$ws.Range("G60").Value2 = 27.44862204242059
$ws.Range("K60").Value2 = 37.41315100629097
$ws.Range("G61").Value2 = 30.34773555612693
$ws.Range("K61").Value2 = 40.77430212015324
$ws.Range("G62").Value2 = 33.83286349733605
$ws.Range("K62").Value2 = 101.2206750385831
